# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-519) from 2023-10-06 (serial 45205) to 2023-10-07 (serial 45206).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C519").Value2 = 45206
